$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Metadata sheet: bump Version and Date
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "2.0.0"

# "2025-09-22" looks like a date, and a plain .Value assignment would get
# auto-converted to a date serial number. Stage it in a throw-away, text
# formatted cell, then copy/paste-values into the real target so the
# destination keeps its own style and receives a literal text value.
$helper = $meta.Range("Z1")
$helper.NumberFormat = "@"
$helper.Value = "2025-09-22"
$helper.Copy()
$meta.Range("B8").PasteSpecial(-4163)
$helper.Clear()

# ---------------------------------------------------------------------------
# 2. Include #1 sheet: add the missing "Description" for the "other" concept
# ---------------------------------------------------------------------------
$include1 = $wb.Worksheets.Item("Include #1")
$include1.Range("B2").Value = "Other encounter class"

# ---------------------------------------------------------------------------
# 3. New "Exclude #2" sheet, modelled on "Include #1" (same layout/styles),
#    placed right after it.
# ---------------------------------------------------------------------------
$include1.Copy($null, $include1)
$exclude2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$exclude2.Name = "Exclude #2"

$exclude2.Range("A2").Value = "OBSENC"
$exclude2.Range("B2").ClearContents()
$exclude2.Range("B4").Value = "http://terminology.hl7.org/CodeSystem/v3-ActCode"

# Restore the original active sheet/tab selection.
$meta.Activate()
